$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.26580536365509
$ws.Range("B1").Value = 2.302929401397705
$ws.Range("C1").Value = 3.806153535842896
$ws.Range("D1").Value = 2.794040441513062
$ws.Range("E1").Value = 1.354339241981506
